# Auto-generated Excel COM-interop script
# Applies the 'Phantom_Profits' pricing-recalculation update described in the commit
# (re-pulled Universalis market-board averages -> new NQ/HQ price & profit figures)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting-leve tables.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 23599
$ws.Range("J3").Value = 23599
$ws.Range("L3").Value = 23599
$ws.Range("N3").Value = -23827

$ws.Range("H5").Value = 171.5
$ws.Range("I5").Value = 214.875
$ws.Range("J5").Value = 84.75
$ws.Range("K5").Value = 214.875
$ws.Range("L5").Value = 84.75
$ws.Range("M5").Value = -99.875
$ws.Range("N5").Value = -314.75

$ws.Range("H21").Value = 30000
$ws.Range("I21").Value = 30000
$ws.Range("K21").Value = 30000
$ws.Range("M21").Value = -29532

$ws.Range("H23").Value = 30000
$ws.Range("I23").Value = 30000
$ws.Range("K23").Value = 30000
$ws.Range("M23").Value = -29766

$ws.Range("H64").Value = 5125
$ws.Range("I64").Value = 4531.1665
$ws.Range("J64").Value = 5837.6
$ws.Range("K64").Value = 4531.1665
$ws.Range("L64").Value = 5837.6
$ws.Range("M64").Value = -4283.1665
$ws.Range("N64").Value = -6333.6

$ws.Range("H67").Value = 5125
$ws.Range("I67").Value = 4531.1665
$ws.Range("J67").Value = 5837.6
$ws.Range("K67").Value = 4531.1665
$ws.Range("L67").Value = 5837.6
$ws.Range("M67").Value = -3673.1665
$ws.Range("N67").Value = -7553.6

$ws.Range("H88").Value = 1666.8334
$ws.Range("J88").Value = 1237.2
$ws.Range("L88").Value = 1237.2
$ws.Range("N88").Value = -2049.2

$ws.Range("H91").Value = 1666.8334
$ws.Range("J91").Value = 1237.2
$ws.Range("L91").Value = 1237.2
$ws.Range("N91").Value = -4045.2

$ws.Range("H102").Value = 23599
$ws.Range("J102").Value = 23599
$ws.Range("L102").Value = 23599
$ws.Range("N102").Value = -30089

$ws.Range("H137").Value = 4294.28
$ws.Range("I137").Value = 3838.7273
$ws.Range("J137").Value = 7635
$ws.Range("K137").Value = 11516.1819
$ws.Range("L137").Value = 22905
$ws.Range("M137").Value = -8966.1819
$ws.Range("N137").Value = -28005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2249.7856
$ws.Range("I2").Value = 2062.125
$ws.Range("K2").Value = 2062.125
$ws.Range("M2").Value = -1949.125

$ws.Range("H41").Value = 23000
$ws.Range("I41").Value = 10000
$ws.Range("K41").Value = 10000
$ws.Range("M41").Value = -9586

$ws.Range("H61").Value = 2839.7576
$ws.Range("I61").Value = 2422.7407
$ws.Range("J61").Value = 4716.3335
$ws.Range("K61").Value = 2422.7407
$ws.Range("L61").Value = 4716.3335
$ws.Range("M61").Value = -2210.7407
$ws.Range("N61").Value = -5140.3335

$ws.Range("H63").Value = 13757.267
$ws.Range("I63").Value = 12204.917
$ws.Range("K63").Value = 12204.917
$ws.Range("M63").Value = -11518.917

$ws.Range("H66").Value = 13757.267
$ws.Range("I66").Value = 12204.917
$ws.Range("K66").Value = 61024.585
$ws.Range("M66").Value = -57592.585

$ws.Range("H101").Value = 23666.334
$ws.Range("J101").Value = 23666.334
$ws.Range("L101").Value = 23666.334
$ws.Range("N101").Value = -30156.334

$ws.Range("H116").Value = 2249.7856
$ws.Range("I116").Value = 2062.125
$ws.Range("K116").Value = 2062.125
$ws.Range("M116").Value = 231.875

$ws.Range("H132").Value = 3293.9666
$ws.Range("I132").Value = 3156.074
$ws.Range("K132").Value = 9468.222
$ws.Range("M132").Value = -6938.222

$ws.Range("H136").Value = 2839.7576
$ws.Range("I136").Value = 2422.7407
$ws.Range("J136").Value = 4716.3335
$ws.Range("K136").Value = 7268.222099999999
$ws.Range("L136").Value = 14149.0005
$ws.Range("M136").Value = -4718.222099999999
$ws.Range("N136").Value = -19249.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2249.7856
$ws.Range("I3").Value = 2062.125
$ws.Range("K3").Value = 2062.125
$ws.Range("M3").Value = -1948.125

$ws.Range("H22").Value = 309.8
$ws.Range("I22").Value = 309.8
$ws.Range("K22").Value = 309.8
$ws.Range("M22").Value = -136.8

$ws.Range("H25").Value = 4406.5
$ws.Range("I25").Value = 4406.5
$ws.Range("K25").Value = 4406.5
$ws.Range("M25").Value = -4171.5

$ws.Range("H86").Value = 5142.0625
$ws.Range("I86").Value = 2733.8572
$ws.Range("K86").Value = 2733.8572
$ws.Range("M86").Value = -1610.8572

$ws.Range("H89").Value = 5142.0625
$ws.Range("I89").Value = 2733.8572
$ws.Range("K89").Value = 13669.286
$ws.Range("M89").Value = -8053.286

$ws.Range("H141").Value = 17495
$ws.Range("I141").Value = 14995
$ws.Range("J141").Value = 19995
$ws.Range("K141").Value = 14995
$ws.Range("L141").Value = 19995
$ws.Range("M141").Value = -9815
$ws.Range("N141").Value = -30355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 179.8077
$ws.Range("I7").Value = 150.25
$ws.Range("J7").Value = 278.33334
$ws.Range("K7").Value = 150.25
$ws.Range("L7").Value = 278.33334
$ws.Range("M7").Value = -37.25
$ws.Range("N7").Value = -504.33334

$ws.Range("H31").Value = 3522.7646
$ws.Range("I31").Value = 2328.8333
$ws.Range("J31").Value = 6388.2
$ws.Range("K31").Value = 2328.8333
$ws.Range("L31").Value = 6388.2
$ws.Range("M31").Value = -2033.8333
$ws.Range("N31").Value = -6978.2

$ws.Range("H34").Value = 3522.7646
$ws.Range("I34").Value = 2328.8333
$ws.Range("J34").Value = 6388.2
$ws.Range("K34").Value = 2328.8333
$ws.Range("L34").Value = 6388.2
$ws.Range("M34").Value = -2126.8333
$ws.Range("N34").Value = -6792.2

$ws.Range("H62").Value = 4594.6665
$ws.Range("I62").Value = 3281
$ws.Range("J62").Value = 6236.75
$ws.Range("K62").Value = 3281
$ws.Range("L62").Value = 6236.75
$ws.Range("M62").Value = -2657
$ws.Range("N62").Value = -7484.75

$ws.Range("H65").Value = 4594.6665
$ws.Range("I65").Value = 3281
$ws.Range("J65").Value = 6236.75
$ws.Range("K65").Value = 16405
$ws.Range("L65").Value = 31183.75
$ws.Range("M65").Value = -13285
$ws.Range("N65").Value = -37423.75

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H135").Value = 79999.664
$ws.Range("J135").Value = 79999
$ws.Range("L135").Value = 79999
$ws.Range("N135").Value = -90139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 125.6
$ws.Range("I2").Value = 122
$ws.Range("K2").Value = 732
$ws.Range("M2").Value = -619

$ws.Range("H126").Value = 6106.1665
$ws.Range("I126").Value = 2427.4
$ws.Range("J126").Value = 24500
$ws.Range("K126").Value = 7282.200000000001
$ws.Range("L126").Value = 73500
$ws.Range("M126").Value = -2342.200000000001
$ws.Range("N126").Value = -83380

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3041.5
$ws.Range("I80").Value = 2750
$ws.Range("K80").Value = 2750
$ws.Range("M80").Value = -1752

$ws.Range("H83").Value = 3041.5
$ws.Range("I83").Value = 2750
$ws.Range("K83").Value = 13750
$ws.Range("M83").Value = -8758

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H55").Value = 727.4286
$ws.Range("I55").Value = 191.75
$ws.Range("J55").Value = 1441.6666
$ws.Range("K55").Value = 191.75
$ws.Range("L55").Value = 1441.6666
$ws.Range("M55").Value = -18.75
$ws.Range("N55").Value = -1787.6666

$ws.Range("H92").Value = 50389
$ws.Range("J92").Value = 50389
$ws.Range("L92").Value = 50389
$ws.Range("N92").Value = -55381

$ws.Range("H100").Value = 963.7
$ws.Range("I100").Value = 925.8889
$ws.Range("K100").Value = 925.8889
$ws.Range("M100").Value = -384.8889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H46").Value = 51665.332
$ws.Range("J46").Value = 51665.332
$ws.Range("L46").Value = 51665.332
$ws.Range("N46").Value = -52127.332

$ws.Range("H134").Value = 51665.332
$ws.Range("J134").Value = 51665.332
$ws.Range("L134").Value = 154995.996
$ws.Range("N134").Value = -160065.996
